# Update the dSF (column F) values for the specific rows that were
# repulled/recalculated, per the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 0
    22 = 2
    29 = -2
    30 = 1
    34 = 3
    37 = 1
    39 = 2
    41 = 0
    45 = -4
    55 = 0
    62 = -1
    63 = 4
    64 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
